# feat: add 2022-Q3 data
#
# 1) "总计" (summary) sheet: insert a new top data row for 2022-Q3 and push the
#    existing 2022-Q2 / 2022-Q1 / 2021-Q4 rows down by one.
# 2) Insert a brand-new worksheet named "2022-Q3" right after "总计" (so the
#    tab order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4) and populate it
#    with the fund-holdings detail for the new quarter.

function Set-TextValue($cell, [string]$text) {
    # Force the value to stay a text cell (avoids Excel's automatic
    # number/leading-zero coercion for numeric-looking strings like "009439"
    # or "4.10"), then drop back to the default "Normal" style so we don't
    # leave a stray quote-prefix style behind on a cell that should carry no
    # explicit style at all.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" sheet — shift existing rows down and insert 2022-Q3 at the top
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Clone the formatting of the existing index column (A2, style-only) into the
# new row 5 before we start moving data, so every data row (2-5) ends up with
# the same bordered/bold "index" style as the original rows.
$summary.Range("A2:D2").Copy()
$summary.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats

# Push data down: row4 -> row5, row3 -> row4, row2 -> row3 (old -> new), then
# write the brand-new 2022-Q3 figures into row 2.
$summary.Cells.Item(5, 1).Value = 3
Set-TextValue $summary.Cells.Item(5, 2) "2021-Q4"
$summary.Cells.Item(5, 3).Value = 3
$summary.Cells.Item(5, 4).Value = 0.06

$summary.Cells.Item(4, 1).Value = 2
Set-TextValue $summary.Cells.Item(4, 2) "2022-Q1"
$summary.Cells.Item(4, 3).Value = 40
$summary.Cells.Item(4, 4).Value = 2.59

$summary.Cells.Item(3, 1).Value = 1
Set-TextValue $summary.Cells.Item(3, 2) "2022-Q2"
$summary.Cells.Item(3, 3).Value = 18
$summary.Cells.Item(3, 4).Value = 1.84

$summary.Cells.Item(2, 1).Value = 0
Set-TextValue $summary.Cells.Item(2, 2) "2022-Q3"
$summary.Cells.Item(2, 3).Value = 9
$summary.Cells.Item(2, 4).Value = 0.42

# ---------------------------------------------------------------------------
# 2. New "2022-Q3" worksheet, inserted right after "总计"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Borrow the header-row and index-column formatting from the existing
# "2022-Q2" sheet so the new sheet matches the others (bold + bordered
# header row, bold + bordered index column).
$template = $wb.Worksheets.Item("2022-Q2")
$template.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$template.Range("A2").Copy()
$q3.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats

Set-TextValue $q3.Cells.Item(1, 2) "基金代码"
Set-TextValue $q3.Cells.Item(1, 3) "基金名称"
Set-TextValue $q3.Cells.Item(1, 4) "基金规模"
Set-TextValue $q3.Cells.Item(1, 5) "股票总仓位"
Set-TextValue $q3.Cells.Item(1, 6) "仓位占比"
Set-TextValue $q3.Cells.Item(1, 7) "持有市值(亿元)"
Set-TextValue $q3.Cells.Item(1, 8) "仓位排名"

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$rows = @(
    @(0, "501059", "西部利得中证国有企业红利指数增强（LOF）A", "4.10", "87.88", "2.40", "0.0984", 8),
    @(1, "009439", "西部利得中证国有企业红利指数增强（LOF）C", "3.44", "87.88", "2.40", "0.0826", 8),
    @(2, "005561", "创金合信中证红利低波动指数A",               "3.32", "94.12", "2.10", "0.0697", 9),
    @(3, "512890", "华泰柏瑞中证红利低波动ETF",                 "2.60", "99.50", "2.23", "0.0580", 8),
    @(4, "005562", "创金合信中证红利低波动指数C",               "2.19", "94.12", "2.10", "0.0460", 9),
    @(5, "006973", "太平睿盈混合A",                             "3.84", "28.79", "0.92", "0.0353", 4),
    @(6, "007669", "太平睿盈混合C",                             "1.04", "28.79", "0.92", "0.0096", 4),
    @(7, "010658", "海富通欣睿混合C",                           "3.22", "20.40", "0.27", "0.0087", 8),
    @(8, "010657", "海富通欣睿混合A",                           "2.71", "20.40", "0.27", "0.0073", 8)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q3.Cells.Item($r, 2) $row[1]
    Set-TextValue $q3.Cells.Item($r, 3) $row[2]
    Set-TextValue $q3.Cells.Item($r, 4) $row[3]
    Set-TextValue $q3.Cells.Item($r, 5) $row[4]
    Set-TextValue $q3.Cells.Item($r, 6) $row[5]
    Set-TextValue $q3.Cells.Item($r, 7) $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Adding a sheet activates it; restore the original active tab ("总计") so
# the workbook-level view state is left exactly as it was before the edit.
$summary.Activate()

Write-Output "2022-Q3 sheet inserted and 总计 rows shifted."
